$d = $word.ActiveDocument

$targets = @(
    "Way too many systems but fortunately most systems same as in previous vessels.",
    "More-less straight forward. Only Panama arrangement not fulfilled due to aft ship design.",
    "I consider material handling quite smooth. Mostly because it was already third vessel.",
    "Block manufacturing timetable and detail design areas could have been more in line to give extra time for design.",
    "In my opinion communication was smooth between different parties.",
    "We got us well employed by changing the hull structure, suppliers and sub-contractors in third vessel."
)

foreach ($t in $targets) {
    $rng = $d.Content
    $rng.Find.Execute($t, $true, $false, $false, $false, $false, $true, 1, $false, "Design: " + $t, 2)
}
